$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 259 (this shifts the existing rows 259-334
# down to 260-335, keeping all of their values intact).
$ws.Rows("259:259").Insert()

# Populate the newly inserted row with the new daily price record.
$ws.Cells.Item(259, 1).Value  = 10
$ws.Cells.Item(259, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(259, 3).Value  = "La Araucanía"
$ws.Cells.Item(259, 4).Value  = 44463
$ws.Cells.Item(259, 5).Value  = 9
$ws.Cells.Item(259, 6).Value  = "Fruta"
$ws.Cells.Item(259, 7).Value  = 100108
$ws.Cells.Item(259, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(259, 9).Value  = 100108006
$ws.Cells.Item(259, 10).Value = "Plátano"
$ws.Cells.Item(259, 11).Value = "Sin especificar"
$ws.Cells.Item(259, 12).Value = "Pintón"
$ws.Cells.Item(259, 13).Value = 900
$ws.Cells.Item(259, 14).Value = 16000
$ws.Cells.Item(259, 15).Value = 17000
$ws.Cells.Item(259, 16).Value = 16556
$ws.Cells.Item(259, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(259, 18).Value = "Ecuador"
$ws.Cells.Item(259, 19).Value = 828
$ws.Cells.Item(259, 20).Value = 20
